$wb = $excel.ActiveWorkbook

# Update the "Last Updated" timestamp on the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 02:55 PM"

# Update 1 Year (%) column F on the Industry Analysis sheet
$ind = $wb.Worksheets.Item("Industry Analysis")
$ind.Range("F2").Value = 18.476
$ind.Range("F3").Value = -7.7404
$ind.Range("F4").Value = 30.7972
$ind.Range("F5").Value = -50.2266
$ind.Range("F6").Value = 61.9649
$ind.Range("F7").Value = -9.1713
$ind.Range("F8").Value = -3.556
$ind.Range("F9").Value = 38.3509
$ind.Range("F10").Value = -6.2497
$ind.Range("F11").Value = 52.6723
$ind.Range("F12").Value = -6.932
$ind.Range("F13").Value = 17.5662
$ind.Range("F14").Value = -35.5106
$ind.Range("F15").Value = 0.6286
$ind.Range("F16").Value = -3.1514
$ind.Range("F17").Value = -20.6354
$ind.Range("F18").Value = -0.0175
$ind.Range("F19").Value = -26.9255
$ind.Range("F20").Value = 44.703
$ind.Range("F21").Value = 10.0506
$ind.Range("F22").Value = 84.6016
$ind.Range("F23").Value = -54.4868
$ind.Range("F24").Value = -12.8122
$ind.Range("F25").Value = -9.182700000000001
$ind.Range("F26").Value = 5.9529
$ind.Range("F27").Value = -33.2998
$ind.Range("F28").Value = -20.4441
$ind.Range("F29").Value = -17.1514
$ind.Range("F30").Value = 24.527
$ind.Range("F31").Value = 57.6193
$ind.Range("F32").Value = -1.527
$ind.Range("F33").Value = -5.2378
$ind.Range("F34").Value = 27.4054
$ind.Range("F35").Value = 6.7961
$ind.Range("F36").Value = -5.6683
$ind.Range("F37").Value = 1.4178
$ind.Range("F38").Value = -22.4272
$ind.Range("F39").Value = 12.3741
$ind.Range("F40").Value = -5.138
$ind.Range("F41").Value = -0.1825
$ind.Range("F42").Value = 23.2483
$ind.Range("F43").Value = 14.456
$ind.Range("F44").Value = -11.1739
$ind.Range("F45").Value = 27.112
$ind.Range("F46").Value = -5.6252
$ind.Range("F47").Value = -36.5148
$ind.Range("F48").Value = -27.8397
$ind.Range("F49").Value = -25.4424
$ind.Range("F50").Value = -49.1173
$ind.Range("F51").Value = -51.065
$ind.Range("F52").Value = -35.4517
$ind.Range("F53").Value = -11.9879
$ind.Range("F54").Value = -3.0992
$ind.Range("F55").Value = -15.3441
$ind.Range("F56").Value = -25.937
$ind.Range("F57").Value = -29.1486
$ind.Range("F58").Value = -6.4093
$ind.Range("F59").Value = -23.3046
$ind.Range("F60").Value = -11.2657
$ind.Range("F61").Value = -9.777699999999999
$ind.Range("F62").Value = -16.0561
$ind.Range("F63").Value = -9.932499999999999
$ind.Range("F64").Value = 51.8767
$ind.Range("F65").Value = -43.5191
$ind.Range("F66").Value = 13.7315
$ind.Range("F67").Value = 12.6111
$ind.Range("F68").Value = 31.7532
$ind.Range("F69").Value = -19.9577
$ind.Range("F70").Value = -12.9642
$ind.Range("F71").Value = 13.2432
$ind.Range("F72").Value = 2.8232
$ind.Range("F73").Value = -9.179
$ind.Range("F74").Value = -14.2931
$ind.Range("F75").Value = 28.3699
$ind.Range("F76").Value = 45.5868
